$d = $word.ActiveDocument

# "Sep 2014" -> "Sept 2014" (abbreviate using the short month form "Sept")
$d.Content.Find.Execute("Sep", $true, $true, $false, $false, $false, $true, 1, $false, "Sept", 2)

# "July" -> "Jul" (abbreviate the full month name to its short form)
$d.Content.Find.Execute("July", $true, $true, $false, $false, $false, $true, 1, $false, "Jul", 2)
